# fixing and implementing new ideas / v.30.2
# Append 3 new ticket rows (187-189) to the bottom of the tickets log.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("2024-05-22", "12:13:36", "Etiquetadora_2",     "-", "-", "-", "-", "12:13:38", "0:00:02"),
    @("2024-05-22", "12:18:13", "Fallo en elevador",  "-", "-", "-", "-", "12:18:30", "0:00:17"),
    @("2024-05-22", "12:43:29", "Ascensor no sube",   "-", "-", "-", "-", "12:43:32", "0:00:03")
)

$startRow = 187
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowValues = $newRows[$i]
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $col = $j + 1
        $text = $rowValues[$j]
        $cell = $ws.Cells.Item($r, $col)
        if ($col -eq 1) {
            # Column A holds dates formatted as plain text (e.g. "2024-05-22").
            # A leading apostrophe forces Excel to keep it as literal text
            # instead of auto-converting it to a date serial number.
            $cell.Formula = "'" + $text
        } else {
            $cell.Value = $text
        }
    }
}
